$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Entrevistador" column (H) entirely - Fecha/Tiempo shift left
$ws.Columns.Item(8).Delete()

# Update the candidate row (row 2) with the new test data
$ws.Range("A2").Value = "TestPrueba"
$ws.Range("B2").Value = "MAR"
$ws.Range("C2").Value = "rioss"
$ws.Range("D2").Value = "PEPETEST@gmail.com"

# Turn the e-mail address into a real hyperlink (Excel auto-applies the
# built-in Hyperlink style to the cell)
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:PEPETEST@gmail.com")

# Re-fit column D now that it holds the longer e-mail address
$ws.Columns.Item(4).AutoFit()
